# Update gh-pages to output generated at 456a3b4
#
# For both the "展览" (Exhibitions) and "全部类型" (All types) sheets:
#   - a new event row is inserted at row 4 (2024-10-03, 南宁·快看漫画动漫游戏嘉年华
#     KKWORLD-mini（取消）), pushing the existing rows 4..N down by one
#   - the "2024 良牙动漫秋季盛典" row (row 3) gets refreshed 想去人数/最低票价
#     counters (F3: 5346 -> 5356, G3: 62 -> 58)

$wb = $excel.ActiveWorkbook

function Insert-NewEventRow($ws, $rowIndex) {
    # Push rowIndex..N down by one, leaving a blank row at rowIndex.
    $ws.Rows.Item($rowIndex).Insert()

    # Column A holds a simple 0-based sequence number (row-1) styled to match
    # the rest of the column (bold, centered, boxed) - reuse the look of the
    # row immediately above rather than hard-coding a style index.
    $aCell = $ws.Cells.Item($rowIndex, 1)
    $aCell.Value2 = $rowIndex - 1
    $above = $ws.Cells.Item($rowIndex - 1, 1)
    $aCell.Font.Bold = $above.Font.Bold
    $aCell.HorizontalAlignment = $above.HorizontalAlignment
    $aCell.VerticalAlignment = $above.VerticalAlignment
    $aCell.Borders.LineStyle = $above.Borders.LineStyle

    # Column B is a literal "yyyy-mm-dd" text label, not a real date - use a
    # leading quote so Excel's input parser keeps it as text, then reset the
    # cell style back to Normal so no stray "Text" number format sticks.
    $bCell = $ws.Cells.Item($rowIndex, 2)
    $bCell.Value2 = "'2024-10-03"
    $bCell.Style = "Normal"

    $ws.Cells.Item($rowIndex, 3).Value2 = "南宁·快看漫画动漫游戏嘉年华 KKWORLD-mini（取消）"
    $ws.Cells.Item($rowIndex, 4).Value2 = "南宁国际会展中心  南宁国际会展中心"
    $ws.Cells.Item($rowIndex, 5).Value2 = "2024.10.03 09:30-10.04 17:30"
    $ws.Cells.Item($rowIndex, 6).Value2 = 389
    $ws.Cells.Item($rowIndex, 7).Value2 = "不可售"
    $ws.Cells.Item($rowIndex, 8).Value2 = "https://show.bilibili.com/platform/detail.html?id=91043"
    $ws.Cells.Item($rowIndex, 9).Value2 = "//i2.hdslb.com/bfs/openplatform/202408/jEAI96Ev1724123680899.jpeg"
}

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Refresh the counters on the "2024 良牙动漫秋季盛典" row (row 3).
    $ws.Range("F3").Value2 = 5356
    $ws.Range("G3").Value2 = 58

    # Insert the new "快看漫画" event as row 4, shifting later rows down.
    Insert-NewEventRow $ws 4
}

Write-Output "done"
